# Updated cryptos list on Mon Jan 22 14:37:46 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '40.738.63'
$c.ClearFormats()
$ws.Cells.Item(2, 5).Value = '  -2.54%  '
# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '2.380.35'
$c.ClearFormats()
$ws.Cells.Item(3, 5).Value = '  -3.91%  '
# Row 4
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.ClearFormats()
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '311.56'
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -2.35%  '
# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '87.75'
$c.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -5.98%  '
# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '0.529'
$c.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -4.37%  '
# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.07%  '
# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.496'
$c.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  -4.32%  '
# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '0.0840'
$c.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -4.85%  '
# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '30.93'
$c.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -7.35%  '
# Row 12
$ws.Cells.Item(12, 5).Value = '  -1.32%  '
# Row 13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '2.737.69'
$c.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -4.22%  '
# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '6.59'
$c.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -4.98%  '
# Row 15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '15.08'
$c.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  -4.00%  '
# Row 16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '2.360.86'
$c.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  -4.41%  '
# Row 17
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '0.765'
$c.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  -4.28%  '
# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '40.598.09'
$c.ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -2.67%  '
# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '0.0₃0914'
$c.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  -4.26%  '
# Row 20
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '6.17'
$c.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -4.85%  '
# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '68.90'
$c.ClearFormats()
# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '10.93'
$c.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -3.69%  '
# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '233.29'
$c.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  -3.82%  '
# Row 24
$ws.Cells.Item(24, 5).Value = '  -4.33%  '
# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.05%  '
# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '1.82'
$c.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -6.84%  '
# Row 27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '23.92'
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -5.05%  '
# Row 28
$ws.Cells.Item(28, 5).Value = '  -1.90%  '
# Row 29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '9.41'
$c.ClearFormats()
$ws.Cells.Item(29, 5).Value = '  -3.79%  '
# Row 30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '34.01'
$c.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  -8.02%  '
# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '153.02'
$c.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -3.02%  '
# Row 32
$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '5.26'
$c.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  -4.75%  '
# Row 33
$ws.Cells.Item(33, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -0.12%  '
# Row 34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '0.0734'
$c.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -4.29%  '
# Row 35
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '2.44'
$c.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -4.96%  '
# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '0.114'
$c.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -2.45%  '
# Row 37
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '16.12'
$c.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -7.72%  '
# Row 38
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '2.79'
$c.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -4.75%  '
# Row 39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '0.0998'
$c.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  -4.40%  '
# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '1.73'
$c.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -8.02%  '
# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '3.88'
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  -3.83%  '
# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -5.31%  '
# Row 43
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '1.964.33'
$c.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  -1.98%  '
# Row 44
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '0.0272'
$c.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -4.78%  '
# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '17.71'
$c.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -7.91%  '
# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '9.65'
$c.ClearFormats()
$ws.Cells.Item(46, 5).Value = '  +1.35%  '
# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '2.75'
$c.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -8.08%  '
# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '2.607.27'
$c.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -4.05%  '
# Row 49
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '93.60'
$c.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -4.50%  '
# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '72.67'
$c.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -5.61%  '
# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '50.85'
$c.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -2.99%  '
